# Apply the two changes described by the commit:
#  1. Bump the cached "datetimeFigureOut" footer date from 11/1/2012 to
#     11/6/2012 everywhere it appears (the slide master + all 11 slide
#     layouts each carry their own cached copy of that field's text).
#  2. Remove the "TextBox 4" shape on slide 1 that held the
#     "http://vk.com/club33848893" link.

$p = $ppt.ActivePresentation

$oldDate = "11/1/2012"
$newDate = "11/6/2012"

# --- 1. Update the Date Placeholder text on the slide master ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq $oldDate) {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

# --- 2. Update the Date Placeholder text on every slide layout ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 3. Remove the vk.com link textbox ("TextBox 4") from slide 1 ---
$slide1 = $p.Slides.Item(1)
for ($si = $slide1.Shapes.Count; $si -ge 1; $si--) {
    $shp = $slide1.Shapes.Item($si)
    if ($shp.Name -eq "TextBox 4") {
        $shp.Delete()
    }
}
